$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'wrestling pants men'
$ws.Cells.Item(2, 1).Value = 'mens sliding pants'
$ws.Cells.Item(3, 1).Value = 'knee pads for baseball'
$ws.Cells.Item(4, 1).Value = 'sliding pad for baseball'
$ws.Cells.Item(5, 1).Value = 'compression pants padded knees'
$ws.Cells.Item(6, 1).Value = 'goalkeeper knee pads'
$ws.Cells.Item(7, 1).Value = 'tights for basketball youth'
$ws.Cells.Item(8, 1).Value = 'honeycomb basketball knee pads'
$ws.Cells.Item(9, 1).Value = 'knee pad shorts'
$ws.Cells.Item(10, 1).Value = 'mens compression capri leggings'
$ws.Cells.Item(11, 1).Value = 'boys compression pants with knee pads'
$ws.Cells.Item(12, 1).Value = 'rodilleras basketball'
$ws.Cells.Item(13, 1).Value = 'honeycomb knee pads basketball'
$ws.Cells.Item(14, 1).Value = 'pad pants'
$ws.Cells.Item(15, 1).Value = 'long basketball knee pads'
$ws.Cells.Item(16, 1).Value = 'men basketball knee'
$ws.Cells.Item(17, 1).Value = 'mens baseball sliding short'
$ws.Cells.Item(18, 1).Value = 'knee pad for basketball'
$ws.Cells.Item(19, 1).Value = 'hex knee pads'
$ws.Cells.Item(20, 1).Value = 'construction pants with knee pads'
$ws.Cells.Item(21, 1).Value = 'men volleyball knee pads'
$ws.Cells.Item(22, 1).Value = 'tight for boys basketball'
$ws.Cells.Item(23, 1).Value = 'bjj kneepads'
$ws.Cells.Item(24, 1).Value = 'boys padded compression pants'
$ws.Cells.Item(25, 1).Value = 'basketball knee tights'
$ws.Cells.Item(26, 1).Value = 'youth basketball compression leggings with knee pads'
$ws.Cells.Item(27, 1).Value = 'basketball knee compression'
$ws.Cells.Item(28, 1).Value = 'youth basketball compression tights'
$ws.Cells.Item(29, 1).Value = 'basketball hex knee pads'
$ws.Cells.Item(30, 1).Value = 'compression basketball leggings'
$ws.Cells.Item(31, 1).Value = 'knee pads for softball'
$ws.Cells.Item(32, 1).Value = 'the best knee pads'
$ws.Cells.Item(33, 1).Value = 'youth basketball knee pad'
$ws.Cells.Item(34, 1).Value = 'mens yoga pants compression'
$ws.Cells.Item(35, 1).Value = 'black knee pads basketball'
$ws.Cells.Item(36, 1).Value = 'boys basketball knee pads youth'
$ws.Cells.Item(37, 1).Value = 'tights basketball'
$ws.Cells.Item(38, 1).Value = 'men capri tights'
$ws.Cells.Item(39, 1).Value = 'nike basketball knee pads'
$ws.Cells.Item(40, 1).Value = 'padded yoga pants'
$ws.Cells.Item(41, 1).Value = 'padded knee pants'
$ws.Cells.Item(42, 1).Value = 'extra padded knee pads'
$ws.Cells.Item(43, 1).Value = 'basketball leg tights'
$ws.Cells.Item(44, 1).Value = 'catchers knee pads'
$ws.Cells.Item(45, 1).Value = 'basketball tights for youth boys'
$ws.Cells.Item(46, 1).Value = 'compression shorts with padding basketball'
$ws.Cells.Item(47, 1).Value = 'wrestling pants youth'
$ws.Cells.Item(48, 1).Value = 'honeycomb knee pads'
$ws.Cells.Item(49, 1).Value = 'baseball leggings for men'
$ws.Cells.Item(50, 1).Value = 'pants men basketball'
$ws.Cells.Item(51, 1).Value = 'long basketball knee pads adult'
$ws.Cells.Item(52, 1).Value = 'yellow leggings for men'
$ws.Cells.Item(53, 1).Value = 'mens compression capri'
$ws.Cells.Item(54, 1).Value = 'wrestling tights youth boy'
$ws.Cells.Item(55, 1).Value = 'basketball youth tights'
$ws.Cells.Item(56, 1).Value = 'need pads for basketball'
$ws.Cells.Item(57, 1).Value = 'compression capris for men'
$ws.Cells.Item(58, 1).Value = 'knee pads for youth basketball'
$ws.Cells.Item(59, 1).Value = 'knee pads under pants'
$ws.Cells.Item(60, 1).Value = 'goalkeeper padded pants'
$ws.Cells.Item(61, 1).Value = 'baseball knee pants'
$ws.Cells.Item(62, 1).Value = 'girls softball leggings'
$ws.Cells.Item(63, 1).Value = 'youth basketball tights'
$ws.Cells.Item(64, 1).Value = 'sliding pad'
$ws.Cells.Item(65, 1).Value = 'work pants with knee pads for men'
$ws.Cells.Item(66, 1).Value = 'multicam combat pants with knee pads'
$ws.Cells.Item(67, 1).Value = 'knee pad honeycomb'
$ws.Cells.Item(68, 1).Value = 'cycling leggings men'
$ws.Cells.Item(69, 1).Value = 'boys basketball compression tights'
$ws.Cells.Item(70, 1).Value = 'compression knee pads work'
$ws.Cells.Item(71, 1).Value = 'softball leggings'
$ws.Cells.Item(72, 1).Value = 'basketball compression tights'
$ws.Cells.Item(73, 1).Value = 'boys basketball pads'
$ws.Cells.Item(74, 1).Value = 'basketball knee pads youth boys black'
$ws.Cells.Item(75, 1).Value = 'lacrosse knee pads'
$ws.Cells.Item(76, 1).Value = 'knee pads youth basketball'
$ws.Cells.Item(77, 1).Value = 'hex kneepads'
$ws.Cells.Item(78, 1).Value = 'knee basketball'
$ws.Cells.Item(79, 1).Value = 'basketball compression pants'
$ws.Cells.Item(80, 1).Value = 'knee pad basketball youth'
$ws.Cells.Item(81, 1).Value = 'weightlifting pants men'
$ws.Cells.Item(82, 1).Value = 'mens knee baseball pants'
$ws.Cells.Item(83, 1).Value = 'thick yoga pad'
$ws.Cells.Item(84, 1).Value = 'padded basketball compression shorts'
$ws.Cells.Item(85, 1).Value = 'under pant knee pads'
$ws.Cells.Item(86, 1).Value = 'youth boys basketball knee pads'
$ws.Cells.Item(87, 1).Value = 'basketball pad'
$ws.Cells.Item(88, 1).Value = 'basketball tights for boys youth'
$ws.Cells.Item(89, 1).Value = 'ready man'
$ws.Cells.Item(90, 1).Value = 'youth basketball knee pads small'
$ws.Cells.Item(91, 1).Value = 'padded compression pants men'
$ws.Cells.Item(92, 1).Value = 'long sliding shorts baseball'
$ws.Cells.Item(93, 1).Value = 'basketball padded compression shorts men'
$ws.Cells.Item(94, 1).Value = 'compression capris men'
$ws.Cells.Item(95, 1).Value = 'gym men leggings'
$ws.Cells.Item(96, 1).Value = 'sliding shorts youth softball'
$ws.Cells.Item(97, 1).Value = 'mens compression pants basketball'
$ws.Cells.Item(98, 1).Value = 'softball tights'
$ws.Cells.Item(99, 1).Value = 'padded basketball shorts'
$ws.Cells.Item(100, 1).Value = 'knee pads basketball men'
